$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 234.875
$ws.Range("I4").Value = 234.875
$ws.Range("K4").Value = 234.875
$ws.Range("M4").Value = -120.875

$ws.Range("H15").Value = 845.5405
$ws.Range("I15").Value = 845.5405
$ws.Range("K15").Value = 2536.6215
$ws.Range("M15").Value = -2367.6215

$ws.Range("H17").Value = 826.5454999999999
$ws.Range("J17").Value = 826.5454999999999
$ws.Range("L17").Value = 2479.6365
$ws.Range("N17").Value = -2815.6365

$ws.Range("H80").Value = 403.66666
$ws.Range("J80").Value = 471.35294
$ws.Range("L80").Value = 1414.05882
$ws.Range("N80").Value = -3410.05882

$ws.Range("H83").Value = 403.66666
$ws.Range("J83").Value = 471.35294
$ws.Range("L83").Value = 4242.17646
$ws.Range("N83").Value = -14226.17646

$ws.Range("H107").Value = 1180.4445
$ws.Range("I107").Value = 1310
$ws.Range("K107").Value = 1310
$ws.Range("M107").Value = 610

$ws.Range("H113").Value = 9102.700000000001
$ws.Range("J113").Value = 8735.799999999999
$ws.Range("L113").Value = 8735.799999999999
$ws.Range("N113").Value = -15243.8

$ws.Range("H132").Value = 44791.707
$ws.Range("I132").Value = 3115.4666
$ws.Range("J132").Value = 114252.11
$ws.Range("K132").Value = 9346.399800000001
$ws.Range("L132").Value = 342756.33
$ws.Range("M132").Value = -6816.399800000001
$ws.Range("N132").Value = -347816.33

$ws.Range("H137").Value = 2013.5217
$ws.Range("I137").Value = 1484.75
$ws.Range("J137").Value = 2295.5334
$ws.Range("K137").Value = 4454.25
$ws.Range("L137").Value = 6886.600199999999
$ws.Range("M137").Value = -1904.25
$ws.Range("N137").Value = -11986.6002

$ws.Range("H138").Value = 6336.2256
$ws.Range("J138").Value = 6569.125
$ws.Range("L138").Value = 19707.375
$ws.Range("N138").Value = -29987.375

$ws.Range("H141").Value = 4725.8096
$ws.Range("I141").Value = 2571.375
$ws.Range("K141").Value = 7714.125
$ws.Range("M141").Value = -2534.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5806.7144
$ws.Range("I32").Value = 3989.24
$ws.Range("J32").Value = 10350.4
$ws.Range("K32").Value = 3989.24
$ws.Range("L32").Value = 10350.4
$ws.Range("M32").Value = -3702.24
$ws.Range("N32").Value = -10924.4

$ws.Range("H45").Value = 1802.8334
$ws.Range("I45").Value = 922.5833
$ws.Range("K45").Value = 922.5833
$ws.Range("M45").Value = -545.5833

$ws.Range("H61").Value = 7704.353
$ws.Range("I61").Value = 9163.666999999999
$ws.Range("K61").Value = 9163.666999999999
$ws.Range("M61").Value = -8951.666999999999

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 7704.353
$ws.Range("I136").Value = 9163.666999999999
$ws.Range("K136").Value = 27491.001
$ws.Range("M136").Value = -24941.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2622.6924
$ws.Range("I107").Value = 2382.125
$ws.Range("J107").Value = 3007.6
$ws.Range("K107").Value = 2382.125
$ws.Range("L107").Value = 3007.6
$ws.Range("M107").Value = -462.125
$ws.Range("N107").Value = -6847.6

$ws.Range("H134").Value = 2545.8
$ws.Range("I134").Value = 2442.2856
$ws.Range("K134").Value = 7326.8568
$ws.Range("M134").Value = -4791.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H31").Value = 2934.1892
$ws.Range("I31").Value = 1078.8334
$ws.Range("K31").Value = 1078.8334
$ws.Range("M31").Value = -783.8334

$ws.Range("H34").Value = 2934.1892
$ws.Range("I34").Value = 1078.8334
$ws.Range("K34").Value = 1078.8334
$ws.Range("M34").Value = -876.8334

$ws.Range("H58").Value = 5872.6
$ws.Range("I58").Value = 5763.467
$ws.Range("K58").Value = 5763.467
$ws.Range("M58").Value = -5560.467

$ws.Range("H132").Value = 3786.6667
$ws.Range("I132").Value = 2569
$ws.Range("J132").Value = 6222
$ws.Range("K132").Value = 7707
$ws.Range("L132").Value = 18666
$ws.Range("M132").Value = -5177
$ws.Range("N132").Value = -23726

$ws.Range("H136").Value = 5872.6
$ws.Range("I136").Value = 5763.467
$ws.Range("K136").Value = 17290.401
$ws.Range("M136").Value = -14740.401

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4781.421
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 4991.5
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 14974.5
$ws.Range("M5").Value = -2888
$ws.Range("N5").Value = -15198.5

$ws.Range("H12").Value = 945.14813
$ws.Range("I12").Value = 556.875
$ws.Range("J12").Value = 1108.6316
$ws.Range("K12").Value = 1670.625
$ws.Range("L12").Value = 3325.8948
$ws.Range("M12").Value = -1497.625
$ws.Range("N12").Value = -3671.8948

$ws.Range("H68").Value = 2135.0527
$ws.Range("J68").Value = 2291.6875
$ws.Range("L68").Value = 6875.0625
$ws.Range("N68").Value = -8497.0625

$ws.Range("H70").Value = 2006
$ws.Range("I70").Value = 2006
$ws.Range("K70").Value = 6018
$ws.Range("M70").Value = -5703

$ws.Range("H71").Value = 2135.0527
$ws.Range("J71").Value = 2291.6875
$ws.Range("L71").Value = 20625.1875
$ws.Range("N71").Value = -28737.1875

$ws.Range("H73").Value = 2006
$ws.Range("I73").Value = 2006
$ws.Range("K73").Value = 6018
$ws.Range("M73").Value = -4926

$ws.Range("H107").Value = 1531.9231
$ws.Range("J107").Value = 1674.6666
$ws.Range("L107").Value = 5023.9998
$ws.Range("N107").Value = -8863.9998

$ws.Range("H135").Value = 4781.421
$ws.Range("I135").Value = 1000
$ws.Range("J135").Value = 4991.5
$ws.Range("K135").Value = 9000
$ws.Range("L135").Value = 44923.5
$ws.Range("M135").Value = -6465
$ws.Range("N135").Value = -49993.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H122").Value = 2446.6453
$ws.Range("I122").Value = 2355.2144
$ws.Range("J122").Value = 2521.9412
$ws.Range("K122").Value = 7065.6432
$ws.Range("L122").Value = 7565.823600000001
$ws.Range("M122").Value = -4615.6432
$ws.Range("N122").Value = -12465.8236

$ws.Range("H126").Value = 3234.6875
$ws.Range("I126").Value = 2605.3
$ws.Range("J126").Value = 4283.6665
$ws.Range("K126").Value = 7815.900000000001
$ws.Range("L126").Value = 12850.9995
$ws.Range("M126").Value = -5345.900000000001
$ws.Range("N126").Value = -17790.9995

$ws.Range("H132").Value = 4122.5
$ws.Range("I132").Value = 3547
$ws.Range("K132").Value = 10641
$ws.Range("M132").Value = -8111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 62000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 62000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 62000
$ws.Range("N42").Value = -63126
$ws.Range("M42").ClearContents()

$ws.Range("H43").Value = 5046875
$ws.Range("J43").Value = 5046875
$ws.Range("L43").Value = 5046875
$ws.Range("N43").Value = -5047261

$ws.Range("H49").Value = 62000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 62000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 62000
$ws.Range("N49").Value = -62294
$ws.Range("M49").ClearContents()

$ws.Range("H132").Value = 4322.364
$ws.Range("I132").Value = 2943.25
$ws.Range("K132").Value = 8829.75
$ws.Range("M132").Value = -6299.75

$ws.Range("H136").Value = 10399.214
$ws.Range("I136").Value = 6871.727
$ws.Range("J136").Value = 23333.334
$ws.Range("K136").Value = 20615.181
$ws.Range("L136").Value = 70000.00199999999
$ws.Range("M136").Value = -18065.181
$ws.Range("N136").Value = -75100.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 4716.5
$ws.Range("I100").Value = 4699.75
$ws.Range("J100").Value = 4750
$ws.Range("K100").Value = 9399.5
$ws.Range("L100").Value = 9500
$ws.Range("M100").Value = -8858.5
$ws.Range("N100").Value = -10582

$ws.Range("H132").Value = 4313.3335
$ws.Range("I132").Value = 4293.5625
$ws.Range("K132").Value = 12880.6875
$ws.Range("M132").Value = -10350.6875

$ws.Range("H136").Value = 4755.857
$ws.Range("I136").Value = 3620.4443
$ws.Range("J136").Value = 6799.6
$ws.Range("K136").Value = 10861.3329
$ws.Range("L136").Value = 20398.8
$ws.Range("M136").Value = -8311.332900000001
$ws.Range("N136").Value = -25498.8
